# Update column F (dSF) values for the "rasmussen_drew" sheet.
# These edits correspond to a repull/recalculation of the dSF data
# (commit message: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = 4
    5  = -4
    6  = -4
    7  = 4
    9  = 4
    10 = 4
    11 = 3
    13 = -1
    14 = 1
    15 = 4
    16 = -3
    18 = 9
    19 = -3
    20 = 2
    21 = -4
    22 = 2
    23 = 4
    24 = 2
    25 = 1
    26 = 5
    27 = 3
    28 = 3
    29 = -3
    30 = 3
    31 = 4
    32 = 1
    33 = -2
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
